# The edit inserts one new data row (a new weekly price observation) right
# above the current row 21, pushing every row from 21..86 down to 22..87.
# Everything else in the sheet (headers, the other 65 existing data rows)
# stays exactly as-is; only the dimension grows from R86 to R87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 21 (and everything below it) down by one row, leaving a clean
# blank row 21 to populate with the new record.
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the new observation's data.
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C21").Value = "Arica y Parinacota"
$ws.Range("D21").Value = 45076
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 100112031
$ws.Range("G21").Value = "Poroto verde"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 1300
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 900
$ws.Range("N21").Value = "$/kilo"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 900
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
